$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.90834366666667
$ws.Range("H2").Value = 95.725031
$ws.Range("I2").Value = 0.1125536485145784
$ws.Range("J2").Value = 0.1157863270269485
$ws.Range("M2").Value = 2.906846333333333
$ws.Range("N2").Value = 8.720538999999999
$ws.Range("O2").Value = 0.005520525738044089
$ws.Range("P2").Value = 0.005624540846623205
$ws.Range("Q2").Value = 92.75265179018987
$ws.Range("R2").Value = 834.7738661117089
$ws.Range("S2").Value = 0.000621355313535498
$ws.Range("T2").Value = 0.0006512449258435446
$ws.Range("G3").Value = 31.90834366666667
$ws.Range("H3").Value = 95.725031
$ws.Range("I3").Value = 0.1125536485145784
$ws.Range("J3").Value = 0.1157863270269485
$ws.Range("O3").Value = 0.3528665483720876
$ws.Range("P3").Value = 0.3595150912979765
$ws.Range("Q3").Value = 5928.657820397802
$ws.Range("R3").Value = 53357.92038358023
$ws.Range("S3").Value = 0.03971641745802444
$ws.Range("T3").Value = 0.04162693193215078
$ws.Range("G4").Value = 31.90834366666667
$ws.Range("H4").Value = 95.725031
$ws.Range("I4").Value = 0.1125536485145784
$ws.Range("J4").Value = 0.1157863270269485
$ws.Range("M4").Value = 137.0717086666666
$ws.Range("N4").Value = 411.2151259999999
$ws.Range("O4").Value = 0.2603191943704447
$ws.Range("P4").Value = 0.2652240042658267
$ws.Range("Q4").Value = 4373.731187113211
$ws.Range("R4").Value = 39363.5806840189
$ws.Range("S4").Value = 0.02929987510476925
$ws.Range("T4").Value = 0.03070931329331981
$ws.Range("G5").Value = 31.90834366666667
$ws.Range("H5").Value = 95.725031
$ws.Range("I5").Value = 0.1125536485145784
$ws.Range("J5").Value = 0.1157863270269485
$ws.Range("M5").Value = 29.2127365
$ws.Range("N5").Value = 58.425473
$ws.Range("O5").Value = 0.05547925319534149
$ws.Range("P5").Value = 0.03768304451958546
$ws.Range("Q5").Value = 932.1300356857771
$ws.Range("R5").Value = 5592.780214114663
$ws.Range("S5").Value = 0.006244392363999767
$ws.Range("T5").Value = 0.004363181316115783
$ws.Range("G6").Value = 31.90834366666667
$ws.Range("H6").Value = 95.725031
$ws.Range("I6").Value = 0.1125536485145784
$ws.Range("J6").Value = 0.1157863270269485
$ws.Range("M6").Value = 171.5584106666666
$ws.Range("N6").Value = 514.6752319999999
$ws.Range("O6").Value = 0.3258144783240821
$ws.Range("P6").Value = 0.331953319069988
$ws.Range("Q6").Value = 5474.144726459132
$ws.Range("R6").Value = 49267.30253813219
$ws.Range("S6").Value = 0.03667160827424947
$ws.Range("T6").Value = 0.03843565555951863
$ws.Range("I7").Value = 0.2312918537506949
$ws.Range("J7").Value = 0.2379348388122522
$ws.Range("M7").Value = 2.906846333333333
$ws.Range("N7").Value = 8.720538999999999
$ws.Range("O7").Value = 0.005520525738044089
$ws.Range("P7").Value = 0.005624540846623205
$ws.Range("Q7").Value = 190.6018423744571
$ws.Range("R7").Value = 1715.416581370114
$ws.Range("S7").Value = 0.00127685263163064
$ws.Range("T7").Value = 0.001338274219734221
$ws.Range("I8").Value = 0.2312918537506949
$ws.Range("J8").Value = 0.2379348388122522
$ws.Range("O8").Value = 0.3528665483720876
$ws.Range("P8").Value = 0.3595150912979765
$ws.Range("S8").Value = 0.08161515809958939
$ws.Range("T8").Value = 0.0855411652985562
$ws.Range("I9").Value = 0.2312918537506949
$ws.Range("J9").Value = 0.2379348388122522
$ws.Range("M9").Value = 137.0717086666666
$ws.Range("N9").Value = 411.2151259999999
$ws.Range("O9").Value = 0.2603191943704447
$ws.Range("P9").Value = 0.2652240042658267
$ws.Range("Q9").Value = 8987.788556171185
$ws.Range("R9").Value = 80890.09700554067
$ws.Range("S9").Value = 0.06020970903282759
$ws.Range("T9").Value = 0.06310603070412958
$ws.Range("I10").Value = 0.2312918537506949
$ws.Range("J10").Value = 0.2379348388122522
$ws.Range("M10").Value = 29.2127365
$ws.Range("N10").Value = 58.425473
$ws.Range("O10").Value = 0.05547925319534149
$ws.Range("P10").Value = 0.03768304451958546
$ws.Range("Q10").Value = 1915.478411724166
$ws.Range("R10").Value = 11492.870470345
$ws.Range("S10").Value = 0.01283189931625469
$ws.Range("T10").Value = 0.00896610912372249
$ws.Range("I11").Value = 0.2312918537506949
$ws.Range("J11").Value = 0.2379348388122522
$ws.Range("M11").Value = 171.5584106666666
$ws.Range("N11").Value = 514.6752319999999
$ws.Range("O11").Value = 0.3258144783240821
$ws.Range("P11").Value = 0.331953319069988
$ws.Range("Q11").Value = 11249.0807556392
$ws.Range("R11").Value = 101241.7268007528
$ws.Range("S11").Value = 0.07535823467039254
$ws.Range("T11").Value = 0.07898325946610975
$ws.Range("G12").Value = 85.57939900000001
$ws.Range("H12").Value = 256.738197
$ws.Range("I12").Value = 0.3018731932863474
$ws.Range("J12").Value = 0.3105433607867011
$ws.Range("M12").Value = 2.906846333333333
$ws.Range("N12").Value = 8.720538999999999
$ws.Range("O12").Value = 0.005520525738044089
$ws.Range("P12").Value = 0.005624540846623205
$ws.Range("Q12").Value = 248.7661621920203
$ws.Range("R12").Value = 2238.895459728183
$ws.Range("S12").Value = 0.001666498733162839
$ws.Range("T12").Value = 0.001746663817392447
$ws.Range("G13").Value = 85.57939900000001
$ws.Range("H13").Value = 256.738197
$ws.Range("I13").Value = 0.3018731932863474
$ws.Range("J13").Value = 0.3105433607867011
$ws.Range("O13").Value = 0.3528665483720876
$ws.Range("P13").Value = 0.3595150912979765
$ws.Range("Q13").Value = 15900.88719259732
$ws.Range("R13").Value = 143107.9847333759
$ws.Range("S13").Value = 0.1065209517610135
$ws.Range("T13").Value = 0.1116450247052113
$ws.Range("G14").Value = 85.57939900000001
$ws.Range("H14").Value = 256.738197
$ws.Range("I14").Value = 0.3018731932863474
$ws.Range("J14").Value = 0.3105433607867011
$ws.Range("M14").Value = 137.0717086666666
$ws.Range("N14").Value = 411.2151259999999
$ws.Range("O14").Value = 0.2603191943704447
$ws.Range("P14").Value = 0.2652240042658267
$ws.Range("Q14").Value = 11730.51444759642
$ws.Range("R14").Value = 105574.6300283678
$ws.Range("S14").Value = 0.07858338647833547
$ws.Range("T14").Value = 0.08236355364601618
$ws.Range("G15").Value = 85.57939900000001
$ws.Range("H15").Value = 256.738197
$ws.Range("I15").Value = 0.3018731932863474
$ws.Range("J15").Value = 0.3105433607867011
$ws.Range("M15").Value = 29.2127365
$ws.Range("N15").Value = 58.425473
$ws.Range("O15").Value = 0.05547925319534149
$ws.Range("P15").Value = 0.03768304451958546
$ws.Range("Q15").Value = 2500.008432815363
$ws.Range("R15").Value = 15000.05059689218
$ws.Range("S15").Value = 0.01674769932321953
$ws.Range("T15").Value = 0.01170221928978695
$ws.Range("G16").Value = 85.57939900000001
$ws.Range("H16").Value = 256.738197
$ws.Range("I16").Value = 0.3018731932863474
$ws.Range("J16").Value = 0.3105433607867011
$ws.Range("M16").Value = 171.5584106666666
$ws.Range("N16").Value = 514.6752319999999
$ws.Range("O16").Value = 0.3258144783240821
$ws.Range("P16").Value = 0.331953319069988
$ws.Range("Q16").Value = 14681.86567824852
$ws.Range("R16").Value = 132136.7911042367
$ws.Range("S16").Value = 0.09835465699061607
$ws.Range("T16").Value = 0.1030858993282942
$ws.Range("G17").Value = 23.7449455
$ws.Range("H17").Value = 47.489891
$ws.Range("I17").Value = 0.08375803763818537
$ws.Range("J17").Value = 0.05744244731349463
$ws.Range("M17").Value = 2.906846333333333
$ws.Range("N17").Value = 8.720538999999999
$ws.Range("O17").Value = 0.005520525738044089
$ws.Range("P17").Value = 0.005624540846623205
$ws.Range("Q17").Value = 69.02290776187482
$ws.Range("R17").Value = 414.1374465712489
$ws.Range("S17").Value = 0.0004623884025496679
$ws.Range("T17").Value = 0.000323087391244752
$ws.Range("G18").Value = 23.7449455
$ws.Range("H18").Value = 47.489891
$ws.Range("I18").Value = 0.08375803763818537
$ws.Range("J18").Value = 0.05744244731349463
$ws.Range("O18").Value = 0.3528665483720876
$ws.Range("P18").Value = 0.3595150912979765
$ws.Range("Q18").Value = 4411.876037945433
$ws.Range("R18").Value = 26471.2562276726
$ws.Range("S18").Value = 0.02955540963980588
$ws.Range("T18").Value = 0.02065142669029023
$ws.Range("G19").Value = 23.7449455
$ws.Range("H19").Value = 47.489891
$ws.Range("I19").Value = 0.08375803763818537
$ws.Range("J19").Value = 0.05744244731349463
$ws.Range("M19").Value = 137.0717086666666
$ws.Range("N19").Value = 411.2151259999999
$ws.Range("O19").Value = 0.2603191943704447
$ws.Range("P19").Value = 0.2652240042658267
$ws.Range("Q19").Value = 3254.760251881877
$ws.Range("R19").Value = 19528.56151129126
$ws.Range("S19").Value = 0.0218038248800218
$ws.Range("T19").Value = 0.01523511589131383
$ws.Range("G20").Value = 23.7449455
$ws.Range("H20").Value = 47.489891
$ws.Range("I20").Value = 0.08375803763818537
$ws.Range("J20").Value = 0.05744244731349463
$ws.Range("M20").Value = 29.2127365
$ws.Range("N20").Value = 58.425473
$ws.Range("O20").Value = 0.05547925319534149
$ws.Range("P20").Value = 0.03768304451958546
$ws.Range("Q20").Value = 693.6548360983608
$ws.Range("R20").Value = 2774.619344393443
$ws.Range("S20").Value = 0.004646833377273828
$ws.Range("T20").Value = 0.00216460629942836
$ws.Range("G21").Value = 23.7449455
$ws.Range("H21").Value = 47.489891
$ws.Range("I21").Value = 0.08375803763818537
$ws.Range("J21").Value = 0.05744244731349463
$ws.Range("M21").Value = 171.5584106666666
$ws.Range("N21").Value = 514.6752319999999
$ws.Range("O21").Value = 0.3258144783240821
$ws.Range("P21").Value = 0.331953319069988
$ws.Range("Q21").Value = 4073.645111346618
$ws.Range("R21").Value = 24441.87066807971
$ws.Range("S21").Value = 0.0272895813385342
$ws.Range("T21").Value = 0.01906821104121746
$ws.Range("G22").Value = 76.69186633333334
$ws.Range("H22").Value = 230.075599
$ws.Range("I22").Value = 0.270523266810194
$ws.Range("J22").Value = 0.2782930260606035
$ws.Range("M22").Value = 2.906846333333333
$ws.Range("N22").Value = 8.720538999999999
$ws.Range("O22").Value = 0.005520525738044089
$ws.Range("P22").Value = 0.005624540846623205
$ws.Range("Q22").Value = 222.9314704475401
$ws.Range("R22").Value = 2006.383234027861
$ws.Range("S22").Value = 0.001493430657165444
$ws.Range("T22").Value = 0.001565270492408241
$ws.Range("G23").Value = 76.69186633333334
$ws.Range("H23").Value = 230.075599
$ws.Range("I23").Value = 0.270523266810194
$ws.Range("J23").Value = 0.2782930260606035
$ws.Range("O23").Value = 0.3528665483720876
$ws.Range("P23").Value = 0.3595150912979765
$ws.Range("Q23").Value = 14249.55923277851
$ws.Range("R23").Value = 128246.0330950066
$ws.Range("S23").Value = 0.0954586114136545
$ws.Range("T23").Value = 0.100050542671768
$ws.Range("G24").Value = 76.69186633333334
$ws.Range("H24").Value = 230.075599
$ws.Range("I24").Value = 0.270523266810194
$ws.Range("J24").Value = 0.2782930260606035
$ws.Range("M24").Value = 137.0717086666666
$ws.Range("N24").Value = 411.2151259999999
$ws.Range("O24").Value = 0.2603191943704447
$ws.Range("P24").Value = 0.2652240042658267
$ws.Range("Q24").Value = 10512.28515914561
$ws.Range("R24").Value = 94610.56643231047
$ws.Range("S24").Value = 0.07042239887449055
$ws.Range("T24").Value = 0.07380999073104733
$ws.Range("G25").Value = 76.69186633333334
$ws.Range("H25").Value = 230.075599
$ws.Range("I25").Value = 0.270523266810194
$ws.Range("J25").Value = 0.2782930260606035
$ws.Range("M25").Value = 29.2127365
$ws.Range("N25").Value = 58.425473
$ws.Range("O25").Value = 0.05547925319534149
$ws.Range("P25").Value = 0.03768304451958546
$ws.Range("Q25").Value = 2240.379282888888
$ws.Range("R25").Value = 13442.27569733333
$ws.Range("S25").Value = 0.01500842881459367
$ws.Range("T25").Value = 0.01048692849053188
$ws.Range("G26").Value = 76.69186633333334
$ws.Range("H26").Value = 230.075599
$ws.Range("I26").Value = 0.270523266810194
$ws.Range("J26").Value = 0.2782930260606035
$ws.Range("M26").Value = 171.5584106666666
$ws.Range("N26").Value = 514.6752319999999
$ws.Range("O26").Value = 0.3258144783240821
$ws.Range("P26").Value = 0.331953319069988
$ws.Range("Q26").Value = 13157.13469920711
$ws.Range("R26").Value = 118414.212292864
$ws.Range("S26").Value = 0.08814039705028984
$ws.Range("T26").Value = 0.09238029367484801
